# Remove the Mansfield product row (and its related hyperlink) from the
# "URL" sheet of Edgesupply.xlsx, since Edgesupply stopped selling
# Mansfield products. The Gerber row (originally row 3) shifts up to
# become row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink attached to the Mansfield row's Link cell (L2)
# before the row shifts, so the stale relationship/hyperlink entry does
# not linger in the worksheet.
$ws.Range("L2").Hyperlinks.Delete()

# Delete the entire Mansfield data row (row 2). The Gerber data that was
# in row 3 moves up to row 2.
$ws.Rows.Item(2).Delete()

# Reset the view: scroll back to the top-left corner (clears the
# topLeftCell="H1" freeze) and move the selection to E13.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E13").Select()
